$wb = $excel.ActiveWorkbook

# --- About sheet: clear stale formatting/empty cells (rows 21-27) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A21:B27").EntireRow.ClearFormats()
$wsAbout.Range("A24").Clear()
$wsAbout.Range("B21:B27").Clear()

# --- Data and Calculations sheet: clear stale custom row format on row 2 ---
$wsData = $wb.Worksheets.Item("Data and Calculations")
$wsData.Range("A2").EntireRow.ClearFormats()

# --- DRC-BDRC: new formulas replacing hardcoded zeros ---
$wsBdrc = $wb.Worksheets.Item("DRC-BDRC")
$wsBdrc.Range("B2").Formula = "='Data and Calculations'!A3"
$wsBdrc.Range("C2").Formula = "=B2"
$wsBdrc.Range("D2:AH2").Formula = "=C2"

# --- DRC-HoDRAUMCUpY: clear stale style on A2 ---
$wsHo = $wb.Worksheets.Item("DRC-HoDRAUMCUpY")
$wsHo.Range("A2").ClearFormats()

# --- DRC-ADRHpDRE: clear stale style on A2 ---
$wsAdr = $wb.Worksheets.Item("DRC-ADRHpDRE")
$wsAdr.Range("A2").ClearFormats()

Write-Host "done"
